# SCD0011 -> SCD0016 update: rename sheet/test-case id and tidy up the
# saved view state, matching the authored diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from SCD0264 to SCD0016 (the file is "SCD0016-038...").
$ws.Name = "SCD0016"

# The TC_ID column (B2:B6) changes from the old "DGS-279" placeholder to
# the real test-case id "SCD0016-038".
$ws.Range("B2:B6").Value = "SCD0016-038"

# Column B is widened (best-fit) to accommodate the longer id text.
$ws.Columns.Item(2).ColumnWidth = 11.5

# Move the live selection to B7, matching the saved sheet view.
$ws.Range("B7").Select()
